$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Consolidate the two separate "GUID" help-file references for the
# Manager Plugins / SDN Plugins rows into a single shared user-guide GUID.
$ws.Range("C17").Value = "GUID-65309889-62B2-43BE-81CE-6A4B650AAFEE"
$ws.Range("C18").Value = "GUID-65309889-62B2-43BE-81CE-6A4B650AAFEE"

# Reflect the author's on-screen selection/scroll position at save time.
$ws.Range("C13").Select()
$excel.ActiveWindow.ScrollRow = 4
